$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(
"2021-10-05 13:38:42.002646", "2021-10-05 13:38:42.002657", "2021-10-05 13:38:42.002660", "2021-10-05 13:38:42.002663", "2021-10-05 13:38:42.002666", "2021-10-05 13:38:42.002668", "2021-10-05 13:38:42.002671", "2021-10-05 13:38:42.002674", "2021-10-05 13:38:42.002677", "2021-10-05 13:38:42.002679", "2021-10-05 13:38:42.002682", "2021-10-05 13:38:42.002684", "2021-10-05 13:38:42.002687", "2021-10-05 13:38:42.002689", "2021-10-05 13:38:42.002692", "2021-10-05 13:38:42.002694", "2021-10-05 13:38:42.002697", "2021-10-05 13:38:42.002700", "2021-10-05 13:38:42.002702", "2021-10-05 13:38:42.002705", "2021-10-05 13:38:42.002707", "2021-10-05 13:38:42.002710", "2021-10-05 13:38:42.002713", "2021-10-05 13:38:42.002715", "2021-10-05 13:38:42.002719", "2021-10-05 13:38:42.002721", "2021-10-05 13:38:42.002724", "2021-10-05 13:38:42.002726", "2021-10-05 13:38:42.002729", "2021-10-05 13:38:42.002732", "2021-10-05 13:38:42.002734", "2021-10-05 13:38:42.002737", "2021-10-05 13:38:42.002740", "2021-10-05 13:38:42.002742", "2021-10-05 13:38:42.002745", "2021-10-05 13:38:42.002747", "2021-10-05 13:38:42.002750", "2021-10-05 13:38:42.002752", "2021-10-05 13:38:42.002755", "2021-10-05 13:38:42.002757", "2021-10-05 13:38:42.002760", "2021-10-05 13:38:42.002763", "2021-10-05 13:38:42.002766", "2021-10-05 13:38:42.002768", "2021-10-05 13:38:42.002771", "2021-10-05 13:38:42.002773", "2021-10-05 13:38:42.002776", "2021-10-05 13:38:42.002779", "2021-10-05 13:38:42.002781", "2021-10-05 13:38:42.002784", "2021-10-05 13:38:42.002786", "2021-10-05 13:38:42.002789", "2021-10-05 13:38:42.002792", "2021-10-05 13:38:42.002795", "2021-10-05 13:38:42.002798", "2021-10-05 13:38:42.002800", "2021-10-05 13:38:42.002803", "2021-10-05 13:38:42.002805", "2021-10-05 13:38:42.002808", "2021-10-05 13:38:42.002810", "2021-10-05 13:38:42.002813", "2021-10-05 13:38:42.002816", "2021-10-05 13:38:42.002818", "2021-10-05 13:38:42.002820", "2021-10-05 13:38:42.002824", "2021-10-05 13:38:42.002827", "2021-10-05 13:38:42.002829", "2021-10-05 13:38:42.002832", "2021-10-05 13:38:42.002835", "2021-10-05 13:38:42.002837", "2021-10-05 13:38:42.002840", "2021-10-05 13:38:42.002843", "2021-10-05 13:38:42.002845", "2021-10-05 13:38:42.002848", "2021-10-05 13:38:42.002850", "2021-10-05 13:38:42.002853", "2021-10-05 13:38:42.002858", "2021-10-05 13:38:42.002861", "2021-10-05 13:38:42.002864", "2021-10-05 13:38:42.002866", "2021-10-05 13:38:42.002869", "2021-10-05 13:38:42.002871", "2021-10-05 13:38:42.002874", "2021-10-05 13:38:42.002877", "2021-10-05 13:38:42.002879", "2021-10-05 13:38:42.002882", "2021-10-05 13:38:42.002884", "2021-10-05 13:38:42.002887", "2021-10-05 13:38:42.002889", "2021-10-05 13:38:42.002892", "2021-10-05 13:38:42.002895", "2021-10-05 13:38:42.002897", "2021-10-05 13:38:42.002901", "2021-10-05 13:38:42.002904", "2021-10-05 13:38:42.002907", "2021-10-05 13:38:42.002909", "2021-10-05 13:38:42.002912", "2021-10-05 13:38:42.002914", "2021-10-05 13:38:42.002917", "2021-10-05 13:38:42.002919", "2021-10-05 13:38:42.002922", "2021-10-05 13:38:42.002924", "2021-10-05 13:38:42.002927", "2021-10-05 13:38:42.002929", "2021-10-05 13:38:42.002932", "2021-10-05 13:38:42.002934", "2021-10-05 13:38:42.002937", "2021-10-05 13:38:42.002940", "2021-10-05 13:38:42.002944", "2021-10-05 13:38:42.002947", "2021-10-05 13:38:42.002950", "2021-10-05 13:38:42.002952", "2021-10-05 13:38:42.002955", "2021-10-05 13:38:42.002957", "2021-10-05 13:38:42.002960", "2021-10-05 13:38:42.002962", "2021-10-05 13:38:42.002965", "2021-10-05 13:38:42.002967", "2021-10-05 13:38:42.002970", "2021-10-05 13:38:42.002972", "2021-10-05 13:38:42.002975", "2021-10-05 13:38:42.002977", "2021-10-05 13:38:42.002980", "2021-10-05 13:38:42.002982", "2021-10-05 13:38:42.002985", "2021-10-05 13:38:42.002987", "2021-10-05 13:38:42.002990", "2021-10-05 13:38:42.002992", "2021-10-05 13:38:42.002997", "2021-10-05 13:38:42.003000", "2021-10-05 13:38:42.003002", "2021-10-05 13:38:42.003005", "2021-10-05 13:38:42.003007", "2021-10-05 13:38:42.003010", "2021-10-05 13:38:42.003012", "2021-10-05 13:38:42.003015", "2021-10-05 13:38:42.003018", "2021-10-05 13:38:42.003020", "2021-10-05 13:38:42.003023", "2021-10-05 13:38:42.003026", "2021-10-05 13:38:42.003028", "2021-10-05 13:38:42.003031", "2021-10-05 13:38:42.003033", "2021-10-05 13:38:42.003036", "2021-10-05 13:38:42.003038", "2021-10-05 13:38:42.003041", "2021-10-05 13:38:42.003043", "2021-10-05 13:38:42.003046", "2021-10-05 13:38:42.003049", "2021-10-05 13:38:42.003052", "2021-10-05 13:38:42.003054", "2021-10-05 13:38:42.003057", "2021-10-05 13:38:42.003059", "2021-10-05 13:38:42.003062", "2021-10-05 13:38:42.003065", "2021-10-05 13:38:42.003067", "2021-10-05 13:38:42.003070", "2021-10-05 13:38:42.003072", "2021-10-05 13:38:42.003075", "2021-10-05 13:38:42.003078", "2021-10-05 13:38:42.003080", "2021-10-05 13:38:42.003083", "2021-10-05 13:38:42.003085", "2021-10-05 13:38:42.003088", "2021-10-05 13:38:42.003091", "2021-10-05 13:38:42.003093", "2021-10-05 13:38:42.003096", "2021-10-05 13:38:42.003098", "2021-10-05 13:38:42.003101", "2021-10-05 13:38:42.003103", "2021-10-05 13:38:42.003106", "2021-10-05 13:38:42.003109", "2021-10-05 13:38:42.003112", "2021-10-05 13:38:42.003116", "2021-10-05 13:38:42.003118", "2021-10-05 13:38:42.003121", "2021-10-05 13:38:42.003123", "2021-10-05 13:38:42.003126", "2021-10-05 13:38:42.003128", "2021-10-05 13:38:42.003131", "2021-10-05 13:38:42.003134", "2021-10-05 13:38:42.003136", "2021-10-05 13:38:42.003139", "2021-10-05 13:38:42.003142", "2021-10-05 13:38:42.003144", "2021-10-05 13:38:42.003147", "2021-10-05 13:38:42.003149", "2021-10-05 13:38:42.003152", "2021-10-05 13:38:42.003155", "2021-10-05 13:38:42.003157", "2021-10-05 13:38:42.003160", "2021-10-05 13:38:42.003162", "2021-10-05 13:38:42.003165", "2021-10-05 13:38:42.003167", "2021-10-05 13:38:42.003170", "2021-10-05 13:38:42.003173", "2021-10-05 13:38:42.003175", "2021-10-05 13:38:42.003178", "2021-10-05 13:38:42.003181", "2021-10-05 13:38:42.003183", "2021-10-05 13:38:42.003186", "2021-10-05 13:38:42.003189", "2021-10-05 13:38:42.003192", "2021-10-05 13:38:42.003194", "2021-10-05 13:38:42.003197", "2021-10-05 13:38:42.003199", "2021-10-05 13:38:42.003202", "2021-10-05 13:38:42.003204", "2021-10-05 13:38:42.003207", "2021-10-05 13:38:42.003210", "2021-10-05 13:38:42.003212", "2021-10-05 13:38:42.003215", "2021-10-05 13:38:42.003217", "2021-10-05 13:38:42.003220", "2021-10-05 13:38:42.003222", "2021-10-05 13:38:42.003225", "2021-10-05 13:38:42.003228", "2021-10-05 13:38:42.003230", "2021-10-05 13:38:42.003233", "2021-10-05 13:38:42.003235", "2021-10-05 13:38:42.003237", "2021-10-05 13:38:42.003240", "2021-10-05 13:38:42.003242", "2021-10-05 13:38:42.003245", "2021-10-05 13:38:42.003247", "2021-10-05 13:38:42.003250", "2021-10-05 13:38:42.003252", "2021-10-05 13:38:42.003255", "2021-10-05 13:38:42.003258", "2021-10-05 13:38:42.003260", "2021-10-05 13:38:42.003263", "2021-10-05 13:38:42.003265", "2021-10-05 13:38:42.003269", "2021-10-05 13:38:42.003272", "2021-10-05 13:38:42.003274", "2021-10-05 13:38:42.003277", "2021-10-05 13:38:42.003279", "2021-10-05 13:38:42.003282", "2021-10-05 13:38:42.003284", "2021-10-05 13:38:42.003287", "2021-10-05 13:38:42.003289", "2021-10-05 13:38:42.003292", "2021-10-05 13:38:42.003294", "2021-10-05 13:38:42.003297", "2021-10-05 13:38:42.003299", "2021-10-05 13:38:42.003302", "2021-10-05 13:38:42.003305", "2021-10-05 13:38:42.003307", "2021-10-05 13:38:42.003310", "2021-10-05 13:38:42.003312", "2021-10-05 13:38:42.003315", "2021-10-05 13:38:42.003317", "2021-10-05 13:38:42.003320", "2021-10-05 13:38:42.003322", "2021-10-05 13:38:42.003325", "2021-10-05 13:38:42.003327", "2021-10-05 13:38:42.003330", "2021-10-05 13:38:42.003332", "2021-10-05 13:38:42.003335", "2021-10-05 13:38:42.003338", "2021-10-05 13:38:42.003340", "2021-10-05 13:38:42.003343"
)

# Header cell F1 - bold/bordered style matching existing headers (copy from E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
